# Fruta / hortaliza, semanal
# Insert a new weekly data row at the top of the "Ciruela" block (row 49),
# shifting the existing rows 49:71 down to 50:72.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("49:49").Insert()

$ws.Range("A49").Value = 5
$ws.Range("B49").Value = "Macroferia Regional de Talca"
$ws.Range("C49").Value = "Maule"
$ws.Range("D49").Value = 44582
$ws.Range("E49").Value = 7
$ws.Range("F49").Value = "Fruta"
$ws.Range("G49").Value = 100103
$ws.Range("H49").Value = "Frutos de hueso (carozo)"
$ws.Range("I49").Value = 100103002
$ws.Range("J49").Value = "Ciruela"
$ws.Range("K49").Value = "Black Amber"
$ws.Range("L49").Value = "Primera"
$ws.Range("M49").Value = 200
$ws.Range("N49").Value = 9000
$ws.Range("O49").Value = 9000
$ws.Range("P49").Value = 9000
$ws.Range("Q49").Value = "$/bandeja 18 kilos granel"
$ws.Range("R49").Value = "Provincia de Curicó"
$ws.Range("S49").Value = 500
$ws.Range("T49").Value = 18
